$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.037708335920549
$ws.Range("D2").Value = 1.046290747192389
$ws.Range("E2").Value = 1.046002545876006
$ws.Range("F2").Value = 1.055882509905798
$ws.Range("I2").Value = 1.035704922036065
$ws.Range("J2").Value = 1.042809738020734
$ws.Range("K2").Value = 1.049056596035517
$ws.Range("L2").Value = 1.048769202523492
$ws.Range("M2").Value = 1.058621746425689
$ws.Range("N2").Value = 1.018040134670662

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.039145625159971
$ws.Range("D3").Value = 1.047009176182324
$ws.Range("E3").Value = 1.047206049672822
$ws.Range("F3").Value = 1.056995609374209
$ws.Range("I3").Value = 1.035837607428554
$ws.Range("J3").Value = 1.043889069687797
$ws.Range("K3").Value = 1.049586880081442
$ws.Range("L3").Value = 1.049783241777533
$ws.Range("M3").Value = 1.059547611738855
$ws.Range("N3").Value = 1.018411740276893

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.040075199399961
$ws.Range("D4").Value = 1.047473187592297
$ws.Range("E4").Value = 1.047984533216033
$ws.Range("F4").Value = 1.057715312902836
$ws.Range("I4").Value = 1.035921775160271
$ws.Range("J4").Value = 1.044586621511811
$ws.Range("K4").Value = 1.049928479879552
$ws.Range("L4").Value = 1.050438563337255
$ws.Range("M4").Value = 1.060145574356796
$ws.Range("N4").Value = 1.018651557742308

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.040465890060417
$ws.Range("D5").Value = 1.047668052115803
$ws.Range("E5").Value = 1.048311747387946
$ws.Range("F5").Value = 1.05801774740347
$ws.Range("I5").Value = 1.035956755007466
$ws.Range("J5").Value = 1.044879672498627
$ws.Range("K5").Value = 1.050071722607669
$ws.Range("L5").Value = 1.050713864519954
$ws.Range("M5").Value = 1.06039668809241
$ws.Range("N5").Value = 1.018752225650712

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.040531482881862
$ws.Range("D6").Value = 1.047700758655648
$ws.Range("E6").Value = 1.048366684633723
$ws.Range("F6").Value = 1.058068519985682
$ws.Range("I6").Value = 1.035962604570767
$ws.Range("J6").Value = 1.044928865425171
$ws.Range("K6").Value = 1.050095752248946
$ws.Range("L6").Value = 1.050760077382194
$ws.Range("M6").Value = 1.060438835426554
$ws.Range("N6").Value = 1.018769119387714

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04008042022376
$ws.Range("D7").Value = 1.047475792190512
$ws.Range("E7").Value = 1.047988905701239
$ws.Range("F7").Value = 1.05771935455143
$ws.Range("I7").Value = 1.03592224415123
$ws.Range("J7").Value = 1.044590538054556
$ws.Range("K7").Value = 1.049930395332981
$ws.Range("L7").Value = 1.050442242692047
$ws.Range("M7").Value = 1.060148930808583
$ws.Range("N7").Value = 1.018652903465361

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.038194170266145
$ws.Range("D8").Value = 1.046533721779266
$ws.Range("E8").Value = 1.046409331987411
$ws.Range("F8").Value = 1.056258801344061
$ws.Range("I8").Value = 1.035750113311614
$ws.Range("J8").Value = 1.0431746807101
$ws.Range("K8").Value = 1.049236124890375
$ws.Range("L8").Value = 1.049112075199224
$ws.Range("M8").Value = 1.058934882844915
$ws.Range("N8").Value = 1.01816585283361

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.034866705822902
$ws.Range("D9").Value = 1.044867081584931
$ws.Range("E9").Value = 1.043623771298484
$ws.Range("F9").Value = 1.053680849799657
$ws.Range("I9").Value = 1.035433864409756
$ws.Range("J9").Value = 1.040673121349413
$ws.Range("K9").Value = 1.048001002553345
$ws.Range("L9").Value = 1.046761678286052
$ws.Range("M9").Value = 1.056786808968533
$ws.Range("N9").Value = 1.017302691691863

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.032645632665048
$ws.Range("D10").Value = 1.04375154724516
$ws.Range("E10").Value = 1.041765102734796
$ws.Range("F10").Value = 1.051959220382117
$ws.Range("I10").Value = 1.035214335586556
$ws.Range("J10").Value = 1.039000749555696
$ws.Range("K10").Value = 1.04716967716346
$ws.Range("L10").Value = 1.045190228947433
$ws.Range("M10").Value = 1.055348751113437
$ws.Range("N10").Value = 1.016723884652698

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.03168315884005
$ws.Range("D11").Value = 1.043267452112489
$ws.Range("E11").Value = 1.040959851217245
$ws.Range("F11").Value = 1.05121299689772
$ws.Range("I11").Value = 1.035117213794324
$ws.Range("J11").Value = 1.038275441834503
$ws.Range("K11").Value = 1.046807821461601
$ws.Range("L11").Value = 1.044508664643433
$ws.Range("M11").Value = 1.054724604635534
$ws.Range("N11").Value = 1.016472442440848

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.031325536830371
$ws.Range("D12").Value = 1.043087477767421
$ws.Range("E12").Value = 1.040660675829982
$ws.Range("F12").Value = 1.050935701531528
$ws.Range("I12").Value = 1.035080828155254
$ws.Range("J12").Value = 1.038005852006006
$ws.Range("K12").Value = 1.046673128084924
$ws.Range("L12").Value = 1.044255330250838
$ws.Range("M12").Value = 1.054492547422187
$ws.Range("N12").Value = 1.016378921928486

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.031402253328812
$ws.Range("D13").Value = 1.043126090115818
$ws.Range("E13").Value = 1.040724853152796
$ws.Range("F13").Value = 1.050995187562552
$ws.Range("I13").Value = 1.035088647041428
$ws.Range("J13").Value = 1.038063688070083
$ws.Range("K13").Value = 1.046702033136365
$ws.Range("L13").Value = 1.044309679132853
$ws.Range("M13").Value = 1.054542334540423
$ws.Range("N13").Value = 1.016398988012654

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.03165360010219
$ws.Range("D14").Value = 1.04325257864039
$ws.Range("E14").Value = 1.040935122717625
$ws.Range("F14").Value = 1.05119007793183
$ws.Range("I14").Value = 1.035114212476482
$ws.Range("J14").Value = 1.038253161114654
$ws.Range("K14").Value = 1.046796693463613
$ws.Range("L14").Value = 1.044487727452552
$ws.Range("M14").Value = 1.054705427244666
$ws.Range("N14").Value = 1.016464714536364

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.031808447542884
$ws.Range("D15").Value = 1.043330491200131
$ws.Range("E15").Value = 1.04106466746765
$ws.Range("F15").Value = 1.051310141039272
$ws.Range("I15").Value = 1.03512992306543
$ws.Range("J15").Value = 1.038369877970953
$ws.Range("K15").Value = 1.046854979155107
$ws.Range("L15").Value = 1.044597406130798
$ws.Range("M15").Value = 1.054805884630411
$ws.Range("N15").Value = 1.016505194390819

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.032709492894223
$ws.Range("D16").Value = 1.043783652648947
$ws.Range("E16").Value = 1.041818535133776
$ws.Range("F16").Value = 1.052008728828337
$ws.Range("I16").Value = 1.035220737734985
$ws.Range("J16").Value = 1.03904886106468
$ws.Range("K16").Value = 1.047193652536095
$ws.Range("L16").Value = 1.04523543832021
$ws.Range("M16").Value = 1.055390142739586
$ws.Range("N16").Value = 1.016740554768322

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.033274493946378
$ws.Range("D17").Value = 1.044067624412008
$ws.Range("E17").Value = 1.042291297165928
$ws.Range("F17").Value = 1.052446732799515
$ws.Range("I17").Value = 1.03527715059634
$ws.Range("J17").Value = 1.039474455908371
$ws.Range("K17").Value = 1.047405587917267
$ws.Range("L17").Value = 1.045635357818706
$ws.Range("M17").Value = 1.055756240407649
$ws.Range("N17").Value = 1.016887971055926

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.033603979042729
$ws.Range("D18").Value = 1.044233157955808
$ws.Range("E18").Value = 1.042567009464592
$ws.Range("F18").Value = 1.052702141246773
$ws.Range("I18").Value = 1.035309856092629
$ws.Range("J18").Value = 1.039722586559208
$ws.Range("K18").Value = 1.047529024282885
$ws.Range("L18").Value = 1.045868516610114
$ws.Range("M18").Value = 1.055969638476743
$ws.Range("N18").Value = 1.016973877936863

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.033716313176554
$ws.Range("D19").Value = 1.04428958326209
$ws.Range("E19").Value = 1.042661013253755
$ws.Range("F19").Value = 1.052789216815508
$ws.Range("I19").Value = 1.035320974042558
$ws.Range("J19").Value = 1.03980717388681
$ws.Range("K19").Value = 1.04757108205005
$ws.Range("L19").Value = 1.045947999601563
$ws.Range("M19").Value = 1.056042377942597
$ws.Range("N19").Value = 1.01700315667788

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.033213882012278
$ws.Range("D20").Value = 1.044037167522087
$ws.Range("E20").Value = 1.04224057862386
$ws.Range("F20").Value = 1.052399746570142
$ws.Range("I20").Value = 1.035271118632452
$ws.Range("J20").Value = 1.039428805184158
$ws.Range("K20").Value = 1.047382868082116
$ws.Range("L20").Value = 1.045592461383459
$ws.Range("M20").Value = 1.055716976143284
$ws.Range("N20").Value = 1.016872162816802

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.031579588003208
$ws.Range("D21").Value = 1.043215335338676
$ws.Range("E21").Value = 1.040873205509089
$ws.Range("F21").Value = 1.051132690763696
$ws.Range("I21").Value = 1.03510669265944
$ws.Range("J21").Value = 1.038197370978345
$ws.Range("K21").Value = 1.046768826203319
$ws.Range("L21").Value = 1.044435301409259
$ws.Range("M21").Value = 1.054657406639787
$ws.Range("N21").Value = 1.016445363137505

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.030551366167626
$ws.Range("D22").Value = 1.042697691726803
$ws.Range("E22").Value = 1.040013082030764
$ws.Range("F22").Value = 1.050335377358292
$ws.Range("I22").Value = 1.035001516263025
$ws.Range("J22").Value = 1.03742208613417
$ws.Range("K22").Value = 1.046381109897767
$ws.Range("L22").Value = 1.043706757616233
$ws.Range("M22").Value = 1.053989930969176
$ws.Range("N22").Value = 1.016176301046767

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.031096512144219
$ws.Range("D23").Value = 1.042972192223993
$ws.Range("E23").Value = 1.040469089063915
$ws.Range("F23").Value = 1.05075811219024
$ws.Range("I23").Value = 1.035057442476712
$ws.Range("J23").Value = 1.037833178587525
$ws.Range("K23").Value = 1.046586801684044
$ws.Range("L23").Value = 1.044093067557715
$ws.Range("M23").Value = 1.054343894795311
$ws.Range("N23").Value = 1.016319004255618

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.033241270154384
$ws.Range("D24").Value = 1.044050929996412
$ws.Range("E24").Value = 1.042263496279438
$ws.Range("F24").Value = 1.052420977847203
$ws.Range("I24").Value = 1.035273844832857
$ws.Range("J24").Value = 1.03944943312667
$ws.Range("K24").Value = 1.047393134760213
$ws.Range("L24").Value = 1.04561184476858
$ws.Range("M24").Value = 1.055734718408797
$ws.Range("N24").Value = 1.016879306122451

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.035727399520845
$ws.Range("D25").Value = 1.045298730043293
$ws.Range("E25").Value = 1.044344179024689
$ws.Range("F25").Value = 1.054347830811061
$ws.Range("I25").Value = 1.035517154989385
$ws.Range("J25").Value = 1.041320641825925
$ws.Range("K25").Value = 1.048321703716982
$ws.Range("L25").Value = 1.047370096284175
$ws.Range("M25").Value = 1.057343188147165
$ws.Range("N25").Value = 1.017526428352252

Write-Output "applied vm_pu update for case with 380 kV"
